$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D="327.05"; E="-1.74%" },
    @{ Row=3; D="44.60"; E="1.89%" },
    @{ Row=4; D="5.516"; E="-5.67%" },
    @{ Row=5; D="0.08088"; E="-3.04%" },
    @{ Row=6; D="8.677"; E="-1.49%" },
    @{ Row=7; D="1.914"; E="-3.73%" },
    @{ Row=8; D="4.288"; E="-4.73%" },
    @{ Row=9; D="2.690"; E="-7.18%" },
    @{ Row=10; D="0.9427"; E="0.59%" },
    @{ Row=11; D="0.1192"; E="-4.11%" },
    @{ Row=12; D="0.1866"; E="-4.30%" },
    @{ Row=13; D="0.09936"; E="2.73%" },
    @{ Row=14; D="0.04259"; E="0.02%" },
    @{ Row=15; D="0.1067"; E="-0.01%" },
    @{ Row=16; D="0.001278"; E="-2.20%" },
    @{ Row=17; B="CoinExToken"; C="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D="0.04192"; E="-5.14%" },
    @{ Row=18; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.006025"; E="0.13%" },
    @{ Row=19; B="HotbitToken"; C="https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; D="0.004538"; E="2.53%" },
    @{ Row=20; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.586"; E="2.61%" },
    @{ Row=21; B="BitpandaEcosystemToken"; C="https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D="0.3495"; E="-0.41%" },
    @{ Row=22; B="MCDex"; C="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D="8.349"; E="-4.96%" },
    @{ Row=23; B="ProBitToken"; C="https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D="0.1371"; E="0.63%" },
    @{ Row=24; B="ZBToken"; C="https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D="0.2526"; E="-4.05%" },
    @{ Row=25; B="BitKan"; C="https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"; D="0.001237"; E="-1.65%" },
    @{ Row=26; D="0.0001180"; E="-0.96%" },
    @{ Row=27; D="0.0003992"; E="-0.07%" },
    @{ Row=28 },
    @{ Row=29 },
    @{ Row=30 },
    @{ Row=31 },
    @{ Row=32 },
    @{ Row=33 },
    @{ Row=34 },
    @{ Row=35 },
    @{ Row=36 },
    @{ Row=37 },
    @{ Row=38 },
    @{ Row=39; D="0.02636"; E="-5.84%" },
    @{ Row=40; D="0.05461"; E="-5.60%" },
    @{ Row=41; D="0.007621"; E="-3.97%" },
    @{ Row=42; D="0.1398"; E="-2.18%" },
    @{ Row=43; D="0.007171"; E="-20.68%" },
    @{ Row=44; D="0.002023"; E="-3.75%" },
    @{ Row=45; D="0.008827"; E="-10.58%" },
    @{ Row=46; D="0.00007106"; E="-2.00%" },
    @{ Row=47; D="0.00000000751"; E="-0.06%" },
    @{ Row=48; D="0.003666"; E="13.16%" },
    @{ Row=49; D="0.002272"; E="-0.36%" },
    @{ Row=50; D="0.00002102"; E="-0.06%" },
    @{ Row=51; D="0.0002002"; E="-0.06%" }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) { $ws.Cells.Item($r.Row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($r.Row, 4).Value = "'" + $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($r.Row, 5).Value = "'" + $r.E }
    $ws.Cells.Item($r.Row, 7).Value = "'8"
}
